$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data to match latest scrape.
# Cells are forced to Text format ("@") before assignment so that
# numeric-looking strings (e.g. "322.50", "0.600") keep their exact
# textual representation instead of being coerced into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.914.99"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.339.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.67"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.07%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.403"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.919.85"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.950.65"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.76%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.308.12"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "434.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.28%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.65"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.74"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.84"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.90%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.84"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.38%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.43"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.49"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.839.65"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.36%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.24"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.42%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.28"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.68"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "322.50"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0274"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.990"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.20"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.21%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.78%  "
